$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column E (in_service) to TRUE for rows 10 through 15
$ws.Range("E10:E15").Value = $true
